$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-10-17"

# Update the header label cell (shared string "2022 (through 10-16)" -> "2022 (through 10-17)")
$ws.Range("I1").Value = "2022 (through 10-17)"

# Update the October (row 10) and November (row 11) 2022 values, and the Total (row 14)
$ws.Range("I10").Value = 146
$ws.Range("I11").Value = 56
$ws.Range("I14").Value = 1334
